# Adds the two new "personel" sheets (MÜHENDİS, TEKNİKER), fills them with
# the new staff names, and removes the leftover fixed row heights on the
# four existing "birim"/neighbourhood sheets (their rows now autosize using
# the sheet's default row height instead of a hard-coded ht="30"/"45"/"60").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Reset any explicit row heights back to the sheet default on the
#    four existing sheets (AKDENİZ, MEZİTLİ, TOROSLAR, YENİŞEHİR).
# ---------------------------------------------------------------------
foreach ($name in @("AKDENİZ", "MEZİTLİ", "TOROSLAR", "YENİŞEHİR")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.UsedRange.Rows.AutoFit() | Out-Null
}

# ---------------------------------------------------------------------
# 2. Add the "MÜHENDİS" sheet (engineers) after the last existing sheet.
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sMuhendis = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$sMuhendis.Name = "MÜHENDİS"

$sMuhendis.Range("A1").Value = "Lokman ALKAN"
$sMuhendis.Range("A3").Value = "Kemal KORKMAZ"
$sMuhendis.Range("A5").Value = "Cihan KARA"
$sMuhendis.Range("A6").Value = "Tevfik YILDIZ"
$sMuhendis.Range("A7").Value = "Sariye KUŞÇU"
$sMuhendis.Range("A8").Value = "Mehmet Nejat AY"
$sMuhendis.Range("A9").Value = "Gülizar YILDIZ"
$sMuhendis.Range("A10").Value = "Mustafa DIKI"
$sMuhendis.Range("A2").Value = "H. Nilgün KIYMAÇ"
$sMuhendis.Range("A4").Value = "Hüseyin KURT"

$sMuhendis.Columns.Item(1).AutoFit() | Out-Null
$sMuhendis.Range("D9").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Add the "TEKNİKER" sheet (technicians) after MÜHENDİS.
# ---------------------------------------------------------------------
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sTekniker = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet2)
$sTekniker.Name = "TEKNİKER"

$sTekniker.Range("A1").Value = "Mustafa BAYSAL"
$sTekniker.Range("A2").Value = "Engin UĞURLU"
$sTekniker.Range("A3").Value = "Özlem AYDINLI"
$sTekniker.Range("A4").Value = "Şekip KORKMAZ"
$sTekniker.Range("A5").Value = "Dudu POYRAZ"
$sTekniker.Range("A6").Value = "Hüseyin POYRAZ"
$sTekniker.Range("A7").Value = "Mehmet CEYLAN"
$sTekniker.Range("A8").Value = "Nurhan ARSLAN"
$sTekniker.Range("A9").Value = "Şehmus ÖZTÜRK"
$sTekniker.Range("A10").Value = "Mustafa GÜRBÜZ"
$sTekniker.Range("A11").Value = "Ahmet ÇELİK"
$sTekniker.Range("A12").Value = "Seval ÇELİK"
$sTekniker.Range("A13").Value = "Müjde TÜRKMEN"
$sTekniker.Range("A14").Value = "Hasan İN"

$sTekniker.Columns.Item(1).AutoFit() | Out-Null
$sTekniker.Range("D10").Select() | Out-Null

# TEKNİKER ends up as the active/selected tab, matching the saved workbook.
$sTekniker.Activate()
